$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 34: DG / NL / 35.037114000000003 / NH3+H2O / POS / [M+NH4]+ / [M+H-H2O]+ / TG ammonium loss + water ----
$ws.Range("B34").Value = "NL"
$ws.Range("C34").Value = 35.037114000000003
$ws.Range("D34").Value = "NH3+H2O"
$ws.Range("E34").Value = "POS"
$ws.Range("F34").Value = "[M+NH4]+"
$ws.Range("H34").Value = "TG ammonium loss + water"
$ws.Range("G34").Value = "[M+H-H2O]+"

# ---- Row 35: DG / NL / 18.010565 / H2O / POS / [M+H]+ / [M+H-H2O]+ / TG water loss ----
$ws.Range("B35").Value = "NL"
$ws.Range("C35").Value = 18.010565
$ws.Range("D35").Value = "H2O"
$ws.Range("E35").Value = "POS"
$ws.Range("F35").Value = "[M+H]+"
$ws.Range("H35").Value = "TG water loss"
$ws.Range("G35").Value = "[M+H-H2O]+"

# ---- CLASS column filled last (introduces the "DG" shared string) ----
$ws.Range("A34").Value = "DG"
$ws.Range("A35").Value = "DG"

# ---- Styling to match the rest of the table ----
# Column C (EXACTMASS) keeps the sheet-wide numeric format, handled first so
# the later blanket alignment pass does not inherit column C's default style.
$massRange = $ws.Range("C34:C35")
$massRange.NumberFormat = "0.0000"
$massRange.HorizontalAlignment = -4152   # xlRight

$ws.Range("A34:B35").HorizontalAlignment = -4108   # xlCenter
$ws.Range("D34:G35").HorizontalAlignment = -4108   # xlCenter

$remarksRange = $ws.Range("H34:H35")
$remarksRange.HorizontalAlignment = -4131  # xlLeft

# ---- View state to match the saved workbook ----
$ws.Range("A35").Select()
